$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Audi A4 B6 1.9 TDI 130 cv Gasóleo 2000 - 2004 AVF, AWX"

$errText = @"
Error: Message: stale element reference: stale element not found
  (Session info: chrome=130.0.6723.70); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#stale-element-reference-exception
Stacktrace:
	GetHandleVerifier [0x00007FF616893AB5+28005]
	(No symbol) [0x00007FF6167F83B0]
	(No symbol) [0x00007FF61669580A]
	(No symbol) [0x00007FF6166A7D0B]
	(No symbol) [0x00007FF61669CF04]
	(No symbol) [0x00007FF61669AE09]
	(No symbol) [0x00007FF61669E5E7]
	(No symbol) [0x00007FF61669E6A0]
	(No symbol) [0x00007FF6166E7B4B]
	(No symbol) [0x00007FF6166DB3AC]
	(No symbol) [0x00007FF61670BA3A]
	(No symbol) [0x00007FF6166D9246]
	(No symbol) [0x00007FF61670BC50]
	(No symbol) [0x00007FF61672B8B3]
	(No symbol) [0x00007FF61670B7E3]
	(No symbol) [0x00007FF6166D75C8]
	(No symbol) [0x00007FF6166D8731]
	GetHandleVerifier [0x00007FF616B8643D+3118829]
	GetHandleVerifier [0x00007FF616BD6C90+3448640]
	GetHandleVerifier [0x00007FF616BCCF0D+3408317]
	GetHandleVerifier [0x00007FF61695A40B+841403]
	(No symbol) [0x00007FF61680340F]
	(No symbol) [0x00007FF6167FF484]
	(No symbol) [0x00007FF6167FF61D]
	(No symbol) [0x00007FF6167EEB79]
	BaseThreadInitThunk [0x00007FFB490D7374+20]
	RtlUserThreadStart [0x00007FFB4A61CC91+33]

"@

$ws.Range("B3").Value = $errText
